# Componente.xlsx update
#   - Swap the SMD rotary encoder (EC11J12-15P30C-SW) for the THT variant
#     (EC11E12-15P30C-SW): part number, description, TME/datasheet links
#     and unit price.
#   - Minor cosmetic follow-up: clear the (inherited) fill on the part
#     number cell, nudge column B's width, and leave the selection where
#     the user ended up (C12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: encoder part -------------------------------------------------
$ws.Range("B7").Value = "EC11E12-15P30C-SW"
$ws.Range("C7").Value = "Codor: incremental; THT; 15imp/rotaţie; două semnale A şi B; 5mA"
$ws.Range("D7").Formula = '=HYPERLINK("https://www.tme.eu/ro/details/ec11e12-15p30c-sw/encodere-incrementale/sr-passives/","TME")'
$ws.Range("E7").Formula = '=HYPERLINK("https://www.tme.eu/Document/44e8c47524c4eb6c460cbc7fca5d0c53/EC11E12-15P30C-SW.pdf","Link")'
$ws.Range("F7").Value = 6.26

# Explicitly clear the cell fill (No Fill) - this is what produces the new
# cell style seen on B7 after the edit.
$ws.Range("B7").Interior.ColorIndex = -4142

# --- Column B got a hair wider after the new text was entered -----------
$ws.Columns.Item(2).ColumnWidth = 19.857142857142858

# --- Leave the selection where the editor ended up -----------------------
$ws.Range("C12").Select() | Out-Null
